$wb = $excel.ActiveWorkbook

# --- profile sheet: insert new column N "pro_usda_soil_order" ---
$wsProfile = $wb.Worksheets.Item("profile")
$wsProfile.Columns.Item(14).Insert()

$wsProfile.Range("N1").Value = "pro_usda_soil_order"

# Set Spodosols cells first, then Inceptisols, so shared-string indices
# come out in the same order as the target workbook.
$wsProfile.Range("N6").Value = "Spodosols"
$wsProfile.Range("N7").Value = "Spodosols"
$wsProfile.Range("N8").Value = "Spodosols"
$wsProfile.Range("N4").Value = "Inceptisols"
$wsProfile.Range("N5").Value = "Inceptisols"

# --- metadata sheet: wrap text on the long bibliographical reference ---
$wsMeta = $wb.Worksheets.Item("metadata")
$wsMeta.Range("M4").WrapText = $true
$wsMeta.Rows.Item(4).RowHeight = 409.6

# --- selection / active-sheet bookkeeping to match the saved view state ---
[void]$wsProfile.Range("O11").Select()
[void]$wsMeta.Range("A4").Select()
